$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.290.20"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "2.647.24"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("D5").Value = "'603.51"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'180.28"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").Value = "'0.175"
$ws.Range("E9").Value = "  +4.85%  "
$ws.Range("D10").Value = "2.644.52"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "'0.360"
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("D13").Value = "'5.06"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  +3.75%  "
$ws.Range("D15").Value = "3.111.91"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "72.210.65"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "'26.69"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "2.646.48"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "'11.97"
$ws.Range("E19").Value = "  +4.71%  "
$ws.Range("D20").Value = "'380.42"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "'7.96"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").Value = "'4.20"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  +10.86%  "
$ws.Range("D24").Value = "'73.29"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "'4.41"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'10.11"
$ws.Range("E27").Value = "  +4.22%  "
$ws.Range("D28").Value = "2.782.73"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "0.0₃0958"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").Value = "'525.21"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").Value = "'8.16"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'165.04"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").Value = "'19.40"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").Value = "'0.114"
$ws.Range("E38").Value = "  -5.52%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'1.41"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "'19.07"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").Value = "'1.86"
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("E42").Value = "  +4.12%  "
$ws.Range("D43").Value = "'5.09"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "'0.335"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").Value = "'39.34"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").Value = "'151.86"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").Value = "'3.73"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D49").Value = "'0.547"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("D50").Value = "'1.71"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("E51").Value = "  -2.10%  "
